# Update column G ("K") values for rows 2-57 with the newly regenerated
# strikeout (K) counts, as described in the commit message:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = 1
    3 = 0
    4 = 3
    5 = 1
    6 = 3
    7 = 1
    8 = 2
    9 = 0
    10 = 2
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 3
    16 = 0
    17 = 0
    18 = 2
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 2
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 4
    29 = 1
    30 = 0
    31 = 4
    32 = 4
    33 = 3
    34 = 3
    35 = 0
    36 = 3
    37 = 3
    38 = 2
    39 = 0
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 4
    48 = 1
    49 = 2
    50 = 2
    51 = 1
    52 = 0
    53 = 0
    54 = 1
    55 = 0
    56 = 0
    57 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
